# LNI-317: Add back the previously removed signature blocks to the end
# of the test Statutory Instrument document.
#
# Appends 6 paragraphs right before the final section break:
#   1. empty "N3"-styled paragraph (numbering switched off)
#   2. empty "SigBlock"-styled paragraph
#   3. "SigBlock" paragraph: <tab> + "Senior Official" (Sig_Signee char style)
#   4. "SigBlock" paragraph: <tab> + "A senior officer of the" (Sig_title char style)
#   5. "SigBlock" paragraph: <tab> + "Department of Agriculture, Environment
#      and Rural Affairs" (Sig_title char style)
#   6. empty "N3"-styled paragraph (numbering switched off)

$d = $word.ActiveDocument

# Collapsed range at the very end of the document's main story, i.e.
# immediately after the last existing paragraph and before the
# section properties.
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)

# Build the new paragraphs as OOXML and drop them in in one shot so the
# paragraph/run skeleton (styles, numPr, tabs) comes out exactly right.
$wordml = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newBlockXml = @"
<w:p $wordml><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr></w:p><w:p $wordml><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr></w:p><w:p $wordml><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="SigSignee"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:t>Senior Official</w:t></w:r></w:p><w:p $wordml><w:pPr><w:pStyle w:val="SigBlock"/><w:rPr><w:rStyle w:val="Sigtitle"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:t>A senior officer of the</w:t></w:r></w:p><w:p $wordml><w:pPr><w:pStyle w:val="SigBlock"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>Department of Agriculture, Environment and Rural Affairs</w:t></w:r></w:p><w:p $wordml><w:pPr><w:pStyle w:val="N3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr></w:p>
"@

[void]$insertionPoint.InsertXML($newBlockXml)

# InsertXML doesn't carry the character style reference onto the inline
# text runs themselves (only onto the paragraph mark run properties), so
# apply it explicitly via the Word object model for the three signature
# text runs - this is what produces <w:rPr><w:rStyle .../></w:rPr> on
# each of them, matching a genuine Word edit.
$paraCount = $d.Paragraphs.Count
$signeePara = $d.Paragraphs.Item($paraCount - 3)   # "Senior Official"
$officerPara = $d.Paragraphs.Item($paraCount - 2)  # "A senior officer of the"
$deptPara = $d.Paragraphs.Item($paraCount - 1)     # "Department of Agriculture..."

$signeeTextRange = $signeePara.Range
$signeeTextRange.SetRange($signeeTextRange.Start + 1, $signeeTextRange.End - 1)
$signeeTextRange.Style = "SigSignee"

$officerTextRange = $officerPara.Range
$officerTextRange.SetRange($officerTextRange.Start + 1, $officerTextRange.End - 1)
$officerTextRange.Style = "Sigtitle"

$deptTextRange = $deptPara.Range
$deptTextRange.SetRange($deptTextRange.Start + 1, $deptTextRange.End - 1)
$deptTextRange.Style = "Sigtitle"
